$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.980.40"
$ws.Range("E2").Value = "  -1.60%  "
$ws.Range("D3").Value = "3.382.43"
$ws.Range("E3").Value = "  -0.74%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'573.58"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").Value = "'137.12"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.379.92"
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("E9").Value = "  -1.43%  "
$ws.Range("E10").Value = "  +1.62%  "
$ws.Range("E11").Value = "  -2.32%  "
$ws.Range("D12").Value = "'0.387"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.56%  "
$ws.Range("D13").Value = "3.959.38"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "'26.18"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.79%  "
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").Value = "3.381.19"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "61.131.31"
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("D19").Value = "'14.02"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").Value = "'9.45"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").Value = "'377.01"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.18%  "
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("D24").Value = "3.526.44"
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("E26").Value = "  -2.91%  "
$ws.Range("D27").Value = "'71.25"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("E28").Value = "  +11.22%  "
$ws.Range("E29").Value = "  +8.29%  "
$ws.Range("D30").Value = "'7.51"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "'23.71"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.57%  "
$ws.Range("E36").Value = "  -3.84%  "
$ws.Range("D37").Value = "'1.55"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.50%  "
$ws.Range("D38").Value = "'6.87"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("D39").Value = "'164.41"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E40").Value = "  -3.58%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "'0.774"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.72"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.60%  "
$ws.Range("D44").Value = "'4.42"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("D45").Value = "'41.66"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("E46").Value = "  -2.85%  "
$ws.Range("D47").Value = "'24.16"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.52%  "
$ws.Range("D48").Value = "2.473.07"
$ws.Range("E48").Value = "  +4.06%  "
$ws.Range("D49").Value = "'23.29"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "'6.79"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.85%  "
$ws.Range("E51").Value = "  +3.46%  "
